# Apply updates to the Orders sheet as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6: IsDeleted flag flips from 0 to 1
$ws.Range("F6").Value = 1

# Row 23: Total_Amount and Order_Info updated
$ws.Range("E23").Value = 184.7562627118645
$ws.Range("G23").Value = "Order 22"

# Row 29: Total_Amount updated
$ws.Range("E29").Value = 28988.14117307856

# New row 48: a freshly added order record
$ws.Range("A48").Value = 47
$ws.Range("B48").Value = 2
$ws.Range("C48").Value = 1
$ws.Range("D48").Value = "2025-03-04 07:04:08"
$ws.Range("E48").Value = 0
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = "uPDATE1203"
